$wb = $excel.ActiveWorkbook

# Sheet "TextFileSequence": add a new "name" column after the existing "type" column (F1 -> G1)
$wsText = $wb.Worksheets.Item("TextFileSequence")
$wsText.Range("G1").Value = "name"

# Sheet "Primer": reorder columns from (name, sequence, id, type) to (sequence, id, type, name)
$wsPrimer = $wb.Worksheets.Item("Primer")
$wsPrimer.Range("A1").Value = "sequence"
$wsPrimer.Range("B1").Value = "id"
$wsPrimer.Range("C1").Value = "type"
$wsPrimer.Range("D1").Value = "name"
